$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PARENT_SITE_ID values in rows 2 and 3 from 1002001 to 1036001.
# Force the Text number format first so the numeric-looking id is kept as
# a (shared) string, matching the original cell type, rather than being
# auto-coerced to a number by the Value assignment.
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "1036001"
$ws.Range("A3").Value = "1036001"
